$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gem Mine Planner")

# 1) Update the base reward-per-block output rate.
$ws.Range("G18").Value = 0.0003

# 2) Add the "Output Per ..." breakdown block in columns J:K (rows 18-23).
$ws.Range("J18").Value = "Output Per Block"
$ws.Range("K18").Formula = "=G18"

$ws.Range("J19").Value = "Output Per Second"
$ws.Range("K19").Formula = "=K18/2"

$ws.Range("J20").Value = "Output Per Minute"
$ws.Range("K20").Formula = "=K19*60"

$ws.Range("J21").Value = "Output Per Hour"
$ws.Range("K21").Formula = "=K20*60"

$ws.Range("J22").Value = "Output Per Day"
$ws.Range("K22").Formula = "=K21*24"

$ws.Range("J23").Value = "Output Per Week"
$ws.Range("K23").Formula = "=K22*7"

# 3) Header label for the new column K in the wallet-count table (row 27),
#    matching the wrap-text-only look of the other header cells.
$ws.Range("K27").Value = "Output Per Day Per Wallet Count"
$ws.Range("K27").WrapText = $true

# 4) Fill in "Output Per Day" x wallet-count for every row of the table.
for ($r = 28; $r -le 45; $r++) {
    $ws.Cells.Item($r, 11).Formula = "=`$K`$22*B$r"
}

# 5) Update the active selection to match the latest edit location.
$ws.Activate() | Out-Null
$ws.Range("I22").Select() | Out-Null
